$d = $word.ActiveDocument
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mo="http://schemas.microsoft.com/office/mac/office/2008/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:mv="urn:schemas-microsoft-com:mac:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 wp14"><w:body><w:p w14:paraId="77E3C2E8" w14:textId="77777777" w:rsidR="009B00F9" w:rsidRDefault="00510DE7"><w:r><w:t>When setting up a new alt.chi reviewing system, you have to update some variables.</w:t></w:r></w:p><w:p w14:paraId="77B93C62" w14:textId="77777777" w:rsidR="00510DE7" w:rsidRDefault="00510DE7"/><w:p w14:paraId="184E706B" w14:textId="77777777" w:rsidR="00510DE7" w:rsidRDefault="00510DE7"><w:r><w:t>In System.php:</w:t></w:r></w:p><w:p w14:paraId="143B733E" w14:textId="77777777" w:rsidR="00510DE7" w:rsidRDefault="00510DE7"/><w:p w14:paraId="47C0D81B" w14:textId="77777777" w:rsidR="00510DE7" w:rsidRDefault="00510DE7"><w:r><w:t>Line 22: $submissionsPath</w:t></w:r></w:p><w:p w14:paraId="4C2F61B1" w14:textId="77777777" w:rsidR="00510DE7" w:rsidRDefault="00510DE7"><w:r><w:t>Line 24: $submissionsUrl</w:t></w:r></w:p><w:p w14:paraId="50C3B3FE" w14:textId="77777777" w:rsidR="00510DE7" w:rsidRDefault="00510DE7"><w:r><w:t>Line 134: $mail-&gt;From</w:t></w:r></w:p><w:p w14:paraId="7EDCF84D" w14:textId="77777777" w:rsidR="00510DE7" w:rsidRDefault="00510DE7"><w:r><w:t>Line 135: $mail-&gt;FromName</w:t></w:r></w:p><w:p w14:paraId="7E2EAACA" w14:textId="77777777" w:rsidR="00510DE7" w:rsidRDefault="00510DE7"><w:r><w:t>Line 169: $mail-&gt;From</w:t></w:r></w:p><w:p w14:paraId="5F5AFB4C" w14:textId="77777777" w:rsidR="00510DE7" w:rsidRDefault="00510DE7"><w:r><w:t>Line 170: $mail-&gt;FromName</w:t></w:r></w:p><w:p w14:paraId="40676C73" w14:textId="77777777" w:rsidR="00510DE7" w:rsidRDefault="00510DE7"><w:r><w:t>Line 180: $mail-&gt;Body</w:t></w:r></w:p><w:p w14:paraId="5AB84973" w14:textId="77777777" w:rsidR="00EB7EB9" w:rsidRDefault="00EB7EB9"/><w:p w14:paraId="6AE69A2A" w14:textId="77777777" w:rsidR="00EB7EB9" w:rsidRDefault="00EB7EB9"><w:r><w:t>For author information, will need to alter smarty template author_dump.tpl</w:t></w:r></w:p><w:p w14:paraId="2BB531D9" w14:textId="77777777" w:rsidR="007C6E03" w:rsidRDefault="007C6E03"/><w:p w14:paraId="61461D95" w14:textId="77777777" w:rsidR="007C6E03" w:rsidRDefault="007C6E03"><w:r><w:t>Downformaintenance.tpl</w:t></w:r></w:p><w:p w14:paraId="59E9CF78" w14:textId="77777777" w:rsidR="007C6E03" w:rsidRDefault="007C6E03"><w:r><w:t>Line 9: contact information</w:t></w:r></w:p><w:p w14:paraId="79C90F79" w14:textId="77777777" w:rsidR="007C6E03" w:rsidRDefault="007C6E03"/><w:p w14:paraId="31638FD7" w14:textId="77777777" w:rsidR="007C6E03" w:rsidRDefault="007C6E03"><w:r><w:t>English_review.tpl</w:t></w:r></w:p><w:p w14:paraId="539200BF" w14:textId="77777777" w:rsidR="007C6E03" w:rsidRDefault="007C6E03"><w:r><w:t>Line 11: contact information</w:t></w:r></w:p><w:p w14:paraId="233390D5" w14:textId="77777777" w:rsidR="007C6E03" w:rsidRDefault="007C6E03"/><w:p w14:paraId="0C47EF02" w14:textId="77777777" w:rsidR="007C6E03" w:rsidRDefault="007C6E03"><w:r><w:t>Index.tpl</w:t></w:r></w:p><w:p w14:paraId="3EA8D8F9" w14:textId="77777777" w:rsidR="007C6E03" w:rsidRDefault="007C6E03"><w:r><w:t>Line 6: title</w:t></w:r></w:p><w:p w14:paraId="6E98999A" w14:textId="77777777" w:rsidR="00371DFC" w:rsidRDefault="00371DFC"/><w:p w14:paraId="6F499E48" w14:textId="77777777" w:rsidR="00371DFC" w:rsidRDefault="00371DFC"><w:r><w:t>Introduction.tpl might need some editing.</w:t></w:r><w:r w:rsidR="004A2AD0"><w:t xml:space="preserve"> This is what users see when they’re not logged in.</w:t></w:r></w:p><w:p w14:paraId="764A1B60" w14:textId="77777777" w:rsidR="002C4484" w:rsidRDefault="002C4484"/><w:p w14:paraId="14092F62" w14:textId="77777777" w:rsidR="002C4484" w:rsidRDefault="002C4484"><w:r><w:t>Portal.tpl will the the file that needs editing the most – that is where you let people know if the system is open for submission, reviewin, etc. or if it’s closed.</w:t></w:r><w:r w:rsidR="004A2AD0"><w:t xml:space="preserve"> This is what users see when they are logged in.</w:t></w:r></w:p><w:p w14:paraId="3F6B7FAB" w14:textId="77777777" w:rsidR="0068594C" w:rsidRDefault="0068594C"/><w:p w14:paraId="1A0D9E63" w14:textId="77777777" w:rsidR="0068594C" w:rsidRDefault="0068594C"><w:r><w:t>In show_submission.tpl will need to adjust lines 31-36 for new author names format.</w:t></w:r></w:p><w:p w14:paraId="1781F896" w14:textId="77777777" w:rsidR="0068594C" w:rsidRDefault="0068594C"/><w:p w14:paraId="702A43F1" w14:textId="77777777" w:rsidR="0068594C" w:rsidRDefault="0068594C"><w:r><w:t>In show_submission_guest.tpl will need to do this on lines 19-22</w:t></w:r></w:p><w:p w14:paraId="207D7042" w14:textId="77777777" w:rsidR="0068594C" w:rsidRDefault="0068594C"/><w:p w14:paraId="074869A0" w14:textId="77777777" w:rsidR="0068594C" w:rsidRDefault="008A7205"><w:r><w:t>Start.tpl looks like it’s obsolete.</w:t></w:r></w:p><w:p w14:paraId="0008C3D8" w14:textId="77777777" w:rsidR="008A7205" w:rsidRDefault="008A7205"/><w:p w14:paraId="6B4DE918" w14:textId="77777777" w:rsidR="008A7205" w:rsidRDefault="008A7205"><w:r><w:t>Submission_confirm.tpl and submission_confirm2.tpl will need editing lines 21-29 for new author names. Check also lines 77-83… not sure what’s going on with “extra names”.</w:t></w:r></w:p><w:p w14:paraId="26F41D55" w14:textId="77777777" w:rsidR="009A56E5" w:rsidRDefault="009A56E5"/><w:p w14:paraId="7E15FE72" w14:textId="77777777" w:rsidR="009A56E5" w:rsidRDefault="009A56E5"><w:r><w:t>Submission.tpl might have to change lines 41-48 for new author names.</w:t></w:r></w:p><w:p w14:paraId="2CCF60F0" w14:textId="77777777" w:rsidR="009A56E5" w:rsidRDefault="009A56E5"/><w:p w14:paraId="2AC6A50A" w14:textId="77777777" w:rsidR="009A56E5" w:rsidRDefault="009A56E5"><w:r><w:t>Submisson2.tpl same, 39-41</w:t></w:r></w:p><w:p w14:paraId="21B63FE9" w14:textId="77777777" w:rsidR="0063466F" w:rsidRDefault="0063466F"/><w:p w14:paraId="67BE87DD" w14:textId="77777777" w:rsidR="0063466F" w:rsidRDefault="0063466F"><w:r><w:t>Update_submission.tpl lines 17-20 for new author names.</w:t></w:r></w:p><w:p w14:paraId="136F084C" w14:textId="77777777" w:rsidR="0063466F" w:rsidRDefault="0063466F"/><w:p w14:paraId="259E0B60" w14:textId="77777777" w:rsidR="0063466F" w:rsidRDefault="0063466F"><w:r><w:t>Verified.tpl – line 6, update to 2013</w:t></w:r></w:p><w:p w14:paraId="3B272B75" w14:textId="77777777" w:rsidR="0063466F" w:rsidRDefault="0063466F"/><w:p w14:paraId="68BAE937" w14:textId="77777777" w:rsidR="0063466F" w:rsidRDefault="0063466F"><w:r><w:t>Verify_conditions.tpl line 24 – update to 2013</w:t></w:r></w:p><w:p w14:paraId="6591C848" w14:textId="77777777" w:rsidR="0063466F" w:rsidRDefault="0063466F"/><w:p w14:paraId="24BD622D" w14:textId="77777777" w:rsidR="0063466F" w:rsidRDefault="0063466F"><w:r><w:t>Verify_notagree.tpl line 16 – update organizer names</w:t></w:r></w:p><w:p w14:paraId="16999710" w14:textId="77777777" w:rsidR="008A0866" w:rsidRDefault="008A0866"/><w:p w14:paraId="1D44F21E" w14:textId="38B8BB2A" w:rsidR="008A0866" w:rsidRDefault="008A0866"><w:r><w:t>Changing the status of reviewing (not open yet, open, no longer open) happens in show_submission lines 215-260</w:t></w:r></w:p><w:p w14:paraId="6D69A8E2" w14:textId="77777777" w:rsidR="00F03A51" w:rsidRDefault="00F03A51"/><w:p w14:paraId="0980C6EA" w14:textId="5303A5ED" w:rsidR="00F03A51" w:rsidRDefault="00F03A51"><w:r><w:t>Email information in login.php lines 57-70</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>In db.php you’ll need to set:</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>$dbuser</w:t></w:r></w:p><w:p><w:r><w:t>$dbpass</w:t></w:r></w:p><w:p><w:r><w:t>$dbname</w:t></w:r></w:p><w:p><w:r><w:t>and any other variables that are specific to your database setup.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:sectPr w:rsidR="00F03A51" w:rsidSect="00536DB1"><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1440" w:right="1800" w:bottom="1440" w:left="1800" w:header="708" w:footer="708" w:gutter="0"/><w:cols w:space="708"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
Write-Host "Paragraphs after:" $d.Paragraphs.Count
